$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text (matches source inlineStr cells)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "27.670.93"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").Value = "1.583.97"
$ws.Range("E3").Value = "  -3.20%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "206.65"
$ws.Range("E5").Value = "  -2.40%  "

$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -3.03%  "

$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").Value = "22.29"
$ws.Range("E8").Value = "  -4.75%  "

$ws.Range("D9").Value = "0.254"
$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").Value = "1.808.31"
$ws.Range("E12").Value = "  -3.24%  "

$ws.Range("D13").Value = "1.569.25"
$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("E14").Value = "  -4.12%  "

$ws.Range("E15").Value = "  -5.82%  "

$ws.Range("D16").Value = "27.630.48"
$ws.Range("E16").Value = "  -1.07%  "

$ws.Range("D17").Value = "63.00"
$ws.Range("E17").Value = "  -3.68%  "

$ws.Range("D18").Value = "218.08"
$ws.Range("E18").Value = "  -4.93%  "

$ws.Range("E19").Value = "  -3.59%  "

$ws.Range("D20").Value = "7.32"
$ws.Range("E20").Value = "  -4.99%  "

$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -5.04%  "

$ws.Range("D23").Value = "9.52"
$ws.Range("E23").Value = "  -5.55%  "

$ws.Range("E24").Value = "  -5.05%  "

$ws.Range("D25").Value = "153.46"
$ws.Range("E25").Value = "  -1.43%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").Value = "6.71"
$ws.Range("E27").Value = "  -2.71%  "

$ws.Range("E28").Value = "  -3.15%  "

$ws.Range("E29").Value = "  -4.30%  "

$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("E31").Value = "  -3.61%  "

$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -5.54%  "

$ws.Range("D33").Value = "1.376.82"
$ws.Range("E33").Value = "  -1.46%  "

$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  -5.40%  "

$ws.Range("E35").Value = "  -5.82%  "

$ws.Range("D36").Value = "0.965"
$ws.Range("E36").Value = "  -4.64%  "

$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("E38").Value = "  -3.16%  "

$ws.Range("D39").Value = "0.541"
$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("D40").Value = "0.819"
$ws.Range("E40").Value = "  -3.70%  "

$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("E42").Value = "  -2.37%  "

$ws.Range("E43").Value = "  -3.40%  "

$ws.Range("D44").Value = "63.69"
$ws.Range("E44").Value = "  -3.51%  "

$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").Value = "5.21"
$ws.Range("E46").Value = "  -4.43%  "

$ws.Range("D47").Value = "1.719.43"
$ws.Range("E47").Value = "  -3.29%  "

$ws.Range("D48").Value = "87.59"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("E49").Value = "  -2.26%  "

$ws.Range("D50").Value = "0.0973"
$ws.Range("E50").Value = "  -5.31%  "

$ws.Range("E51").Value = "  -1.44%  "
